$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 0.6848072082125213
$ws.Range("E2").Value = 0.6848072082125213

# Row 3
$ws.Range("D3").Value = 0.01612535653683595
$ws.Range("E3").Value = 0.01612535653683595

# Row 4
$ws.Range("D4").Value = 0.961336995562868
$ws.Range("E4").Value = 0.961336995562868

# Row 5
$ws.Range("D5").Value = 0.0005118484100462969
$ws.Range("E5").Value = 0.0005118484100462969

# Row 6
$ws.Range("D6").Value = 0.2749186407203411
$ws.Range("E6").Value = 0.2749186407203411

# Row 7
$ws.Range("D7").Value = 0.704769956107055
$ws.Range("E7").Value = 0.295230043892945

# Row 8
$ws.Range("D8").Value = 0.9999888348168271
$ws.Range("E8").Value = [double]"1.116518317290094E-05"

# Row 9
$ws.Range("D9").Value = 0.9798358955062828
$ws.Range("E9").Value = 0.02016410449371719

# Row 10
$ws.Range("D10").Value = 0.9999999998040126
$ws.Range("E10").Value = [double]"1.95987448492474E-10"

# Row 11
$ws.Range("D11").Value = 0.9999622489972166
$ws.Range("E11").Value = [double]"3.775100278335941E-05"
$ws.Range("F11").Value = 0.5115985870361328

# Row 12
$ws.Range("D12").Value = 0.8812832669563396
$ws.Range("E12").Value = 0.8812832669563396

# Row 13
$ws.Range("D13").Value = 0.001392269867387415
$ws.Range("E13").Value = 0.001392269867387415

# Row 14
$ws.Range("D14").Value = 0.9922877528399561
$ws.Range("E14").Value = 0.9922877528399561

# Row 15
$ws.Range("D15").Value = 0.0001515278738051976
$ws.Range("E15").Value = 0.0001515278738051976

# Row 16
$ws.Range("D16").Value = 0.07286610765347773
$ws.Range("E16").Value = 0.07286610765347773

# Row 17
$ws.Range("D17").Value = 0.8837349422351977
$ws.Range("E17").Value = 0.1162650577648023

# Row 18
$ws.Range("D18").Value = 0.9999999995212407
$ws.Range("E18").Value = [double]"4.78759254463057E-10"

# Row 19
$ws.Range("D19").Value = 0.9995232383001091
$ws.Range("E19").Value = 0.0004767616998908597

# Row 20
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 0

# Row 21
$ws.Range("D21").Value = 0.9999999999994542
$ws.Range("E21").Value = [double]"5.45785638905727E-13"
$ws.Range("F21").Value = 0.7197238206863403
